$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-06 Friday" "2025-06-07 Saturday"

Replace-Text "957×2=1914" "161×4=644"
Replace-Text "660×6=3960" "643×4=2572"
Replace-Text "759×5=3795" "658×7=4606"
Replace-Text "556×3=1668" "570×7=3990"
Replace-Text "784×4=3136" "156×6=936"
Replace-Text "861×5=4305" "779×9=7011"
Replace-Text "109×4=436" "655×6=3930"
Replace-Text "351×4=1404" "859×5=4295"
Replace-Text "725×3=2175" "905×2=1810"
Replace-Text "750×6=4500" "744×6=4464"
Replace-Text "756×3=2268" "299×8=2392"
Replace-Text "924×4=3696" "933×2=1866"
Replace-Text "837×6=5022" "249×6=1494"
Replace-Text "322×8=2576" "665×5=3325"
Replace-Text "234×5=1170" "440×7=3080"
Replace-Text "313×5=1565" "515×7=3605"
Replace-Text "916×8=7328" "486×7=3402"
Replace-Text "280×2=560" "975×2=1950"
Replace-Text "558×2=1116" "545×5=2725"
Replace-Text "436×8=3488" "233×6=1398"
Replace-Text "300×9=2700" "528×4=2112"
Replace-Text "860×2=1720" "776×2=1552"
Replace-Text "782×7=5474" "938×4=3752"
Replace-Text "910×5=4550" "122×8=976"
Replace-Text "272×6=1632" "114×7=798"

Write-Host "Done"
